$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.247.26'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.591.74'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''212.43'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '''19.03'
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("D11").Value = '''0.0851'
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").Value = '1.815.15'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '1.593.40'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("D16").Value = '''63.86'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '26.239.20'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '''215.77'
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("D20").Value = '''7.32'
$ws.Range("E20").Value = '  -3.04%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '''9.05'
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("D25").Value = '''144.56'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").Value = '''15.12'
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("E30").Value = '  -2.83%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("D33").Value = '1.410.25'
$ws.Range("E33").Value = '  +6.01%  '
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("D36").Value = '''1.47'
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("E37").Value = '  -3.82%  '
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").Value = '''5.81'
$ws.Range("E40").Value = '  +1.87%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").Value = '''0.975'
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").Value = '''0.765'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '1.727.30'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '''60.90'
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").Value = '''86.30'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").Value = '''1.49'
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("E50").Value = '  -2.46%  '
$ws.Range("D51").Value = '''0.999'
$ws.Range("E51").Value = '  -0.08%  '
